$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 77; existing rows 77-222 shift down to 78-223
$ws.Rows.Item(77).Insert()

# Populate the newly inserted row 77 with the new data record
$ws.Cells.Item(77, 1).Value = 4
$ws.Cells.Item(77, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(77, 3).Value = "Los Lagos"
$ws.Cells.Item(77, 4).Value = 44581
$ws.Cells.Item(77, 5).Value = 10
$ws.Cells.Item(77, 6).Value = 100112037
$ws.Cells.Item(77, 7).Value = "Cebollín"
$ws.Cells.Item(77, 8).Value = "Sin especificar"
$ws.Cells.Item(77, 9).Value = "Primera"
$ws.Cells.Item(77, 10).Value = 60
$ws.Cells.Item(77, 11).Value = 6000
$ws.Cells.Item(77, 12).Value = 6000
$ws.Cells.Item(77, 13).Value = 6000
$ws.Cells.Item(77, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(77, 15).Value = "Región Metropolitana"
$ws.Cells.Item(77, 16).Value = 167
$ws.Cells.Item(77, 17).Value = 36
$ws.Cells.Item(77, 18).Value = "Hortaliza"

# Apply the same date number format (index used by column D) to the new D77 cell
$ws.Cells.Item(77, 4).NumberFormat = $ws.Cells.Item(78, 4).NumberFormat
